$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "mark" column (C) with header + three values, matching the
# existing Name/Id columns (A/B) already on the sheet.
$ws.Range("C1").Value = "mark"
$ws.Range("C2").Value = 20
$ws.Range("C3").Value = 80
$ws.Range("C4").Value = 30

# Mirror the author's final selection/active cell from the diff.
$ws.Range("C4").Select()
